$wb = $excel.ActiveWorkbook

# Insert a new "Battery" worksheet right before the existing "Switch" sheet
$switch = $wb.Worksheets.Item("Switch")
$battery = $wb.Worksheets.Add($switch)
$battery.Name = "Battery"

# Header row
$battery.Range("A1").Value = "Battery1_E"
$battery.Range("B1").Value = "Battery1_P"

# Five data rows (7000 / 500), mirroring the other case sheets' shape
for ($r = 2; $r -le 6; $r++) {
    $battery.Cells.Item($r, 1).Value = 7000
    $battery.Cells.Item($r, 2).Value = 500
}

# Update the Reservoir sheet's stored selection
$reservoir = $wb.Worksheets.Item("Reservoir")
[void]$reservoir.Range("A2:A6").Select()

# Battery is the sheet left selected/active after the edits
[void]$battery.Range("C5").Select()
